$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix header typo: "threshhold" -> "threshold"
$ws.Range("E1").Value = "threshold"

# Append 40 new rows (82-121) of categorical summary stats for the "15, 30" threshold group
$data = New-Object 'object[,]' 40,9
$data[0,0]="channelAB";$data[0,1]=0.7;$data[0,2]=1234;$data[0,3]=3;$data[0,4]="15, 30";$data[0,5]=1;$data[0,6]=0.58333333333333304;$data[0,7]=0.35920036204419398;$data[0,8]=0.25
$data[1,0]="channelAB";$data[1,1]=0.7;$data[1,2]=1234;$data[1,3]=3;$data[1,4]="15, 30";$data[1,5]=2;$data[1,6]=0.59523809523809501;$data[1,7]=0.35196828851129802;$data[1,8]=0.45833333333333298
$data[2,0]="channelAB";$data[2,1]=0.7;$data[2,2]=1234;$data[2,3]=3;$data[2,4]="15, 30";$data[2,5]=3;$data[2,6]=0.51428571428571401;$data[2,7]=0.40329853998743398;$data[2,8]=0.41666666666666602
$data[3,0]="channelAB";$data[3,1]=0.7;$data[3,2]=1234;$data[3,3]=3;$data[3,4]="15, 30";$data[3,5]=4;$data[3,6]=0.48;$data[3,7]=0.41937749732716301;$data[3,8]=0.5
$data[4,0]="channelAB";$data[4,1]=0.7;$data[4,2]=1234;$data[4,3]=3;$data[4,4]="15, 30";$data[4,5]=5;$data[4,6]=0.61285714285714199;$data[4,7]=0.368190734559739;$data[4,8]=0.375
$data[5,0]="channelAB";$data[5,1]=0.7;$data[5,2]=1234;$data[5,3]=3;$data[5,4]="15, 30";$data[5,5]=6;$data[5,6]=0.56666666666666599;$data[5,7]=0.36663896644976302;$data[5,8]=0.45833333333333298
$data[6,0]="channelAB";$data[6,1]=0.7;$data[6,2]=1234;$data[6,3]=3;$data[6,4]="15, 30";$data[6,5]=7;$data[6,6]=0.53095238095238095;$data[6,7]=0.42054538283457898;$data[6,8]=0.5
$data[7,0]="channelAB";$data[7,1]=0.7;$data[7,2]=1234;$data[7,3]=3;$data[7,4]="15, 30";$data[7,5]=8;$data[7,6]=0.59809523809523801;$data[7,7]=0.36231760915803701;$data[7,8]=0.41666666666666602
$data[8,0]="channelAB";$data[8,1]=0.7;$data[8,2]=1234;$data[8,3]=3;$data[8,4]="15, 30";$data[8,5]=9;$data[8,6]=0.51619047619047598;$data[8,7]=0.373460047833166;$data[8,8]=0.375
$data[9,0]="channelAB";$data[9,1]=0.7;$data[9,2]=1234;$data[9,3]=3;$data[9,4]="15, 30";$data[9,5]=10;$data[9,6]=0.48380952380952302;$data[9,7]=0.39312151855162197;$data[9,8]=0.58333333333333304
$data[10,0]="channelAB";$data[10,1]=0.8;$data[10,2]=1234;$data[10,3]=3;$data[10,4]="15, 30";$data[10,5]=1;$data[10,6]=0.60595238095238002;$data[10,7]=0.360257777881972;$data[10,8]=0.28571428571428498
$data[11,0]="channelAB";$data[11,1]=0.8;$data[11,2]=1234;$data[11,3]=3;$data[11,4]="15, 30";$data[11,5]=2;$data[11,6]=0.52222222222222203;$data[11,7]=0.38530449680178303;$data[11,8]=0.57142857142857095
$data[12,0]="channelAB";$data[12,1]=0.8;$data[12,2]=1234;$data[12,3]=3;$data[12,4]="15, 30";$data[12,5]=3;$data[12,6]=0.58750000000000002;$data[12,7]=0.39497721450861201;$data[12,8]=0.57142857142857095
$data[13,0]="channelAB";$data[13,1]=0.8;$data[13,2]=1234;$data[13,3]=3;$data[13,4]="15, 30";$data[13,5]=4;$data[13,6]=0.59285714285714197;$data[13,7]=0.37065369935043802;$data[13,8]=0.42857142857142799
$data[14,0]="channelAB";$data[14,1]=0.8;$data[14,2]=1234;$data[14,3]=3;$data[14,4]="15, 30";$data[14,5]=5;$data[14,6]=0.525595238095238;$data[14,7]=0.39785222849909302;$data[14,8]=0.5
$data[15,0]="channelAB";$data[15,1]=0.8;$data[15,2]=1234;$data[15,3]=3;$data[15,4]="15, 30";$data[15,5]=6;$data[15,6]=0.46071428571428502;$data[15,7]=0.38028707332199402;$data[15,8]=0.57142857142857095
$data[16,0]="channelAB";$data[16,1]=0.8;$data[16,2]=1234;$data[16,3]=3;$data[16,4]="15, 30";$data[16,5]=7;$data[16,6]=0.48690476190476101;$data[16,7]=0.39852929217598898;$data[16,8]=0.64285714285714202
$data[17,0]="channelAB";$data[17,1]=0.8;$data[17,2]=1234;$data[17,3]=3;$data[17,4]="15, 30";$data[17,5]=8;$data[17,6]=0.55000000000000004;$data[17,7]=0.40308755491190901;$data[17,8]=0.64285714285714202
$data[18,0]="channelAB";$data[18,1]=0.8;$data[18,2]=1234;$data[18,3]=3;$data[18,4]="15, 30";$data[18,5]=9;$data[18,6]=0.54087301587301495;$data[18,7]=0.40507618263550599;$data[18,8]=0.57142857142857095
$data[19,0]="channelAB";$data[19,1]=0.8;$data[19,2]=1234;$data[19,3]=3;$data[19,4]="15, 30";$data[19,5]=10;$data[19,6]=0.46666666666666601;$data[19,7]=0.417760934199435;$data[19,8]=0.57142857142857095
$data[20,0]="channelAB";$data[20,1]=0.7;$data[20,2]=1234;$data[20,3]=2;$data[20,4]="15, 30";$data[20,5]=1;$data[20,6]=0.79666666666666597;$data[20,7]=0.15564609705922899;$data[20,8]=0.73333333333333295
$data[21,0]="channelAB";$data[21,1]=0.7;$data[21,2]=1234;$data[21,3]=2;$data[21,4]="15, 30";$data[21,5]=2;$data[21,6]=0.71833333333333305;$data[21,7]=0.16902623806633399;$data[21,8]=0.66666666666666596
$data[22,0]="channelAB";$data[22,1]=0.7;$data[22,2]=1234;$data[22,3]=2;$data[22,4]="15, 30";$data[22,5]=3;$data[22,6]=0.63500000000000001;$data[22,7]=0.23366942709202099;$data[22,8]=0.73333333333333295
$data[23,0]="channelAB";$data[23,1]=0.7;$data[23,2]=1234;$data[23,3]=2;$data[23,4]="15, 30";$data[23,5]=4;$data[23,6]=0.68166666666666598;$data[23,7]=0.24509181212450601;$data[23,8]=0.86666666666666603
$data[24,0]="channelAB";$data[24,1]=0.7;$data[24,2]=1234;$data[24,3]=2;$data[24,4]="15, 30";$data[24,5]=5;$data[24,6]=0.82499999999999996;$data[24,7]=0.15643755966498701;$data[24,8]=0.6
$data[25,0]="channelAB";$data[25,1]=0.7;$data[25,2]=1234;$data[25,3]=2;$data[25,4]="15, 30";$data[25,5]=6;$data[25,6]=0.83;$data[25,7]=0.14982702234190701;$data[25,8]=0.53333333333333299
$data[26,0]="channelAB";$data[26,1]=0.7;$data[26,2]=1234;$data[26,3]=2;$data[26,4]="15, 30";$data[26,5]=7;$data[26,6]=0.71833333333333305;$data[26,7]=0.197329652703748;$data[26,8]=0.73333333333333295
$data[27,0]="channelAB";$data[27,1]=0.7;$data[27,2]=1234;$data[27,3]=2;$data[27,4]="15, 30";$data[27,5]=8;$data[27,6]=0.76666666666666605;$data[27,7]=0.16529929171822999;$data[27,8]=0.6
$data[28,0]="channelAB";$data[28,1]=0.7;$data[28,2]=1234;$data[28,3]=2;$data[28,4]="15, 30";$data[28,5]=9;$data[28,6]=0.65666666666666595;$data[28,7]=0.21598115424838901;$data[28,8]=0.66666666666666596
$data[29,0]="channelAB";$data[29,1]=0.7;$data[29,2]=1234;$data[29,3]=2;$data[29,4]="15, 30";$data[29,5]=10;$data[29,6]=0.82166666666666599;$data[29,7]=0.14332613750804199;$data[29,8]=0.6
$data[30,0]="channelAB";$data[30,1]=0.8;$data[30,2]=1234;$data[30,3]=2;$data[30,4]="15, 30";$data[30,5]=1;$data[30,6]=0.82499999999999996;$data[30,7]=0.16115501520387801;$data[30,8]=0.55555555555555503
$data[31,0]="channelAB";$data[31,1]=0.8;$data[31,2]=1234;$data[31,3]=2;$data[31,4]="15, 30";$data[31,5]=2;$data[31,6]=0.69499999999999995;$data[31,7]=0.18429886354876901;$data[31,8]=0.77777777777777701
$data[32,0]="channelAB";$data[32,1]=0.8;$data[32,2]=1234;$data[32,3]=2;$data[32,4]="15, 30";$data[32,5]=3;$data[32,6]=0.788333333333333;$data[32,7]=0.15816544173355901;$data[32,8]=0.55555555555555503
$data[33,0]="channelAB";$data[33,1]=0.8;$data[33,2]=1234;$data[33,3]=2;$data[33,4]="15, 30";$data[33,5]=4;$data[33,6]=0.73499999999999999;$data[33,7]=0.18628141973692899;$data[33,8]=0.66666666666666596
$data[34,0]="channelAB";$data[34,1]=0.8;$data[34,2]=1234;$data[34,3]=2;$data[34,4]="15, 30";$data[34,5]=5;$data[34,6]=0.69833333333333303;$data[34,7]=0.19849029259957299;$data[34,8]=0.77777777777777701
$data[35,0]="channelAB";$data[35,1]=0.8;$data[35,2]=1234;$data[35,3]=2;$data[35,4]="15, 30";$data[35,5]=6;$data[35,6]=0.79;$data[35,7]=0.171673746108428;$data[35,8]=0.55555555555555503
$data[36,0]="channelAB";$data[36,1]=0.8;$data[36,2]=1234;$data[36,3]=2;$data[36,4]="15, 30";$data[36,5]=7;$data[36,6]=0.73;$data[36,7]=0.18660867754823099;$data[36,8]=0.77777777777777701
$data[37,0]="channelAB";$data[37,1]=0.8;$data[37,2]=1234;$data[37,3]=2;$data[37,4]="15, 30";$data[37,5]=8;$data[37,6]=0.71;$data[37,7]=0.205384611843589;$data[37,8]=0.88888888888888795
$data[38,0]="channelAB";$data[38,1]=0.8;$data[38,2]=1234;$data[38,3]=2;$data[38,4]="15, 30";$data[38,5]=9;$data[38,6]=0.72499999999999998;$data[38,7]=0.184019130568876;$data[38,8]=0.77777777777777701
$data[39,0]="channelAB";$data[39,1]=0.8;$data[39,2]=1234;$data[39,3]=2;$data[39,4]="15, 30";$data[39,5]=10;$data[39,6]=0.65833333333333299;$data[39,7]=0.20893617351920099;$data[39,8]=1
$ws.Range("A82:I121").Value = $data

# Update view state to match the final selection
$ws.Range("G6").Select()
